# Apply the "created class diagrams for the first release" edit to the
# Scrum burndown workbook.
#
# The change logs actual effort hours against several tasks on the
# "Sprint Backlog 1 (2)" sheet (Week 2 columns K = Mon, L = Tues), which in
# turn ripples through the per-task totals (column P) and the "Remaining"
# burndown row (row 21) that the burndown chart (chart2) is plotted from.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog 1 (2)")

# --- Log actual hours worked (Week 2: K = Mon, L = Tues) -------------------

# Row 6  - "Task Viewing - list render" : 2 hours on Tuesday
$ws.Range("L6").Value = 2

# Rows 9-13 - Deliverable: UML class diagrams (Task/Team creation, render,
# invite) worked on Monday
$ws.Range("K9").Value = 2
$ws.Range("K10").Value = 2
$ws.Range("K11").Value = 2
$ws.Range("K12").Value = 1
$ws.Range("K13").Value = 1

# Row 16 - "Task Creation" : 6 hours on Tuesday
$ws.Range("L16").Value = 6

# --- Recalculate so the burndown / totals pick up the new actuals ---------
$excel.CalculateFull()

# --- Update the active selection left on the sheet -------------------------
$null = $ws.Range("R10").Select()
